# Insert a new data row at row 144 (shifts existing rows 144:282 down to 145:283)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(144).Insert()

# Fill the new row 144 with the new weekly price record (same constant columns
# as the rest of the dataset for this market/category: Vega Modelo de Temuco / Jengibre)
$ws.Cells.Item(144, 1).Value2 = 10
$ws.Cells.Item(144, 2).Value2 = "Vega Modelo de Temuco"
$ws.Cells.Item(144, 3).Value2 = "La Araucanía"
$ws.Cells.Item(144, 4).Value2 = 45040
$ws.Cells.Item(144, 5).Value2 = 9
$ws.Cells.Item(144, 6).Value2 = 100114007
$ws.Cells.Item(144, 7).Value2 = "Jengibre"
$ws.Cells.Item(144, 8).Value2 = "Sin especificar"
$ws.Cells.Item(144, 9).Value2 = "Primera"
$ws.Cells.Item(144, 10).Value2 = 170
$ws.Cells.Item(144, 11).Value2 = 22000
$ws.Cells.Item(144, 12).Value2 = 25000
$ws.Cells.Item(144, 13).Value2 = 23412
$ws.Cells.Item(144, 14).Value2 = "$/caja 13 kilos"
$ws.Cells.Item(144, 15).Value2 = "Perú"
$ws.Cells.Item(144, 16).Value2 = 1801
$ws.Cells.Item(144, 17).Value2 = 13
$ws.Cells.Item(144, 18).Value2 = "Hortaliza"
